$d = $word.ActiveDocument

# 1. "...du projet Cookie Factory en..." -> "...du projet The Cookie Factory en..."
$d.Content.Find.Execute("projet Cookie Factory en", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "projet The Cookie Factory en", 2)

# 2. "Nous ne recompilons et ne retestons que" -> "Nous ne recompilons et testons que"
$d.Content.Find.Execute("ne retestons", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "testons", 2)

# 3. "le serveur dotNet" -> "le serveur .NET" (both occurrences)
$d.Content.Find.Execute("serveur dotNet", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "serveur .NET", 2)

# 4. "retélécharger" -> "télécharger"
$d.Content.Find.Execute("retélécharger", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "télécharger", 2)
